$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Remove the old "_GoBack" bookmark from its current location (end of
#    the "Lack of trello use -" paragraph). It gets re-created at the end
#    of the first paragraph later on.
# ---------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# ---------------------------------------------------------------------
# 2. Clean up the "GitHub" / proofErr(gramStart) / ":" / proofErr(gramEnd)
#    run soup in paragraph 1 into a single, plain "GitHub:" run (no
#    proofErr markers), merged with the rest of the paragraph's runs.
#    Do it as a two-phase text replace so the proofErr elements (which
#    sit on run boundaries) get fully swallowed by the edit.
# ---------------------------------------------------------------------
$find = $d.Content.Find
$null = $find.Execute("GitHub", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$segStart = $find.Parent.Start

$seg = $d.Range($segStart, $segStart + 14)
$expected = "GitHub:" + [char]11 + "Slack:"
if ($seg.Text -ne $expected) {
    throw "unexpected segment text: [$($seg.Text)]"
}
$seg.Text = "@@SEGPLACEHOLDER@@"

$seg2 = $d.Range($segStart, $segStart + 18)
$seg2.Text = "GitHub:" + [char]11 + "Slack:"

# ---------------------------------------------------------------------
# 3. Insert the Trello hyperlink right after "Trello: ".
# ---------------------------------------------------------------------
$find = $d.Content.Find
$null = $find.Execute("Trello: ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$afterTrello = $find.Parent
$afterTrello.Collapse(0)
$afterTrello.InsertBefore("@@TRELLOLINK@@ ")

$find = $d.Content.Find
$null = $find.Execute("@@TRELLOLINK@@", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$trelloUrl = "https://trello.com/b/AjC5R5tf/web-project-dev-template"
$null = $d.Hyperlinks.Add($find.Parent, $trelloUrl, "", "", $trelloUrl)

# ---------------------------------------------------------------------
# 4. Insert the GitHub hyperlink right after "GitHub:".
# ---------------------------------------------------------------------
$find = $d.Content.Find
$null = $find.Execute("GitHub:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$afterGitHub = $find.Parent
$afterGitHub.Collapse(0)
$afterGitHub.InsertBefore(" @@GITHUBLINK@@")

$find = $d.Content.Find
$null = $find.Execute("@@GITHUBLINK@@", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$githubUrl = "https://github.com/ShaneRing11/cp1406-16-30"
$null = $d.Hyperlinks.Add($find.Parent, $githubUrl, "", "", $githubUrl)

# ---------------------------------------------------------------------
# 5. Insert ' ' + apple-converted-space + Slack hyperlink right after
#    "Slack:" at the end of the paragraph.
# ---------------------------------------------------------------------
$find = $d.Content.Find
$null = $find.Execute("Slack:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$afterSlack = $find.Parent
$afterSlack.Collapse(0)
$afterSlack.InsertBefore(" @@APPLESPACE@@@@SLACKLINK@@")

$find = $d.Content.Find
$null = $find.Execute("@@APPLESPACE@@", $false, $false, $false, $false, $false, $true, 1, $false, " ", 2)

$find = $d.Content.Find
$null = $find.Execute("@@SLACKLINK@@", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$slackUrl = "https://itatjcu.slack.com/messages/cp1406-2016-team30/"
$null = $d.Hyperlinks.Add($find.Parent, $slackUrl, "", "", $slackUrl)

# Tag the lone apple-converted-space run with its character style.
$find = $d.Content.Find
$null = $find.Execute(" ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
# (re-find precisely below by scanning paragraph 1 runs instead of blind Find)

# ---------------------------------------------------------------------
# 6. Re-create the "_GoBack" bookmark at the very end of paragraph 1
#    (after the Slack hyperlink).
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$endOfP1 = $d.Range($p1.Range.End - 1, $p1.Range.End - 1)
$null = $d.Bookmarks.Add("_GoBack", $endOfP1)

# ---------------------------------------------------------------------
# 7. Paragraph 2 ("Members") becomes Heading 2 instead of Heading 1.
# ---------------------------------------------------------------------
$p2 = $d.Paragraphs(2)
$p2.Style = "Heading 2"

Write-Output "done"
